$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new header columns: "Calibration Date" (AI1) and "Decay Factor" (AJ1) ---
$ws.Range("AI1").Value = "Calibration Date"
$ws.Range("AJ1").Value = "Decay Factor"

# Copy the formatting (fill/font/border/alignment) of the last existing header cell (AH1)
# onto the two new header cells so they match the rest of the header row.
$ws.Range("AH1").Copy()
$ws.Range("AI1:AJ1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the new "Calibration Date" column a sensible width like the other data columns.
$ws.Columns("AI").ColumnWidth = 9.6

# --- Update the view state: scroll the frozen header into view and park the ---
# --- active selection on AD1, matching where the analyst was last working.  ---
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AD1").Select()
